# Add barcode alignment test data (commit: "Add barcode alignment to tests")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
# Column A gets wider to fit the new header text; columns B:D become new
# 15-wide columns for the 3-column alignment table.
$ws.Columns("A").ColumnWidth = 41.3
$ws.Range("B1:D1").ColumnWidth = 14.15

# --- Row 16: bold section header ---
$ws.Range("A16").Value = "Horizontal and vertical alignment"
$ws.Range("A16").Font.Bold = $true

# --- Rows 17-21: barcode text sample cells ---
# Writing "12345678" through .Value turns it into a number (Excel's COM
# layer auto-detects numeric-looking strings). To keep it as literal text
# (matching the source data, which stores it as a shared string) we build
# it once via a TEXT() formula in a scratch cell, then paste-special just
# the *values* into the target cells - this carries over the text type
# without leaving behind any NumberFormat/style residue.
$ws.Range("Z1").Formula = "=TEXT(12345678,""0"")"
$ws.Range("Z1").Copy()
$ws.Range("A17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("A18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("A19").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("A21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("C21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("Z1").Clear()

# --- Row 20 / D21: highlighted (orange) empty cells ---
$ws.Range("A20").Interior.Color = 10275833
$ws.Range("D21").Interior.Color = 10275833
